$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C4").Value = "A15289126"
$ws.Range("C4").Select() | Out-Null
